$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.300877666666666
$ws.Range("H2").Value = 12.902633
$ws.Range("I2").Value = 0.04061703229494078
$ws.Range("J2").Value = 0.04061703229494078
$ws.Range("M2").Value = 51.402972
$ws.Range("N2").Value = 154.208916
$ws.Range("O2").Value = 0.2478201393026385
$ws.Range("P2").Value = 0.2478201393026385
$ws.Range("Q2").Value = 221.0778942750919
$ws.Range("R2").Value = 1989.701048475828
$ws.Range("S2").Value = 0.01006571860139199
$ws.Range("T2").Value = 0.01006571860139199

# Row 3
$ws.Range("G3").Value = 4.300877666666666
$ws.Range("H3").Value = 12.902633
$ws.Range("I3").Value = 0.04061703229494078
$ws.Range("J3").Value = 0.04061703229494078
$ws.Range("O3").Value = 0.2505855219821376
$ws.Range("P3").Value = 0.2505855219821376
$ws.Range("Q3").Value = 223.544864802059
$ws.Range("R3").Value = 2011.903783218531
$ws.Range("S3").Value = 0.01017804023899307
$ws.Range("T3").Value = 0.01017804023899308

# Row 4
$ws.Range("G4").Value = 4.300877666666666
$ws.Range("H4").Value = 12.902633
$ws.Range("I4").Value = 0.04061703229494078
$ws.Range("J4").Value = 0.04061703229494078
$ws.Range("M4").Value = 59.18509700000001
$ws.Range("N4").Value = 177.555291
$ws.Range("O4").Value = 0.2853387345614992
$ws.Range("P4").Value = 0.2853387345614991
$ws.Range("Q4").Value = 254.5478618868003
$ws.Range("R4").Value = 2290.930756981203
$ws.Range("S4").Value = 0.01158961259668195
$ws.Range("T4").Value = 0.01158961259668195

# Row 5
$ws.Range("G5").Value = 4.300877666666666
$ws.Range("H5").Value = 12.902633
$ws.Range("I5").Value = 0.04061703229494078
$ws.Range("J5").Value = 0.04061703229494078
$ws.Range("M5").Value = 8.716382666666666
$ws.Range("N5").Value = 26.149148
$ws.Range("O5").Value = 0.04202276799614693
$ws.Range("P5").Value = 0.04202276799614692
$ws.Range("Q5").Value = 37.4880955451871
$ws.Range("R5").Value = 337.3928599066839
$ws.Range("S5").Value = 0.001706840124822304
$ws.Range("T5").Value = 0.001706840124822304

# Row 6
$ws.Range("G6").Value = 4.300877666666666
$ws.Range("H6").Value = 12.902633
$ws.Range("I6").Value = 0.04061703229494078
$ws.Range("J6").Value = 0.04061703229494078
$ws.Range("M6").Value = 36.13945833333333
$ws.Range("N6").Value = 108.418375
$ws.Range("O6").Value = 0.1742328361575779
$ws.Range("P6").Value = 0.1742328361575779
$ws.Range("Q6").Value = 155.4313892312638
$ws.Range("R6").Value = 1398.882503081375
$ws.Range("S6").Value = 0.007076820733051468
$ws.Range("T6").Value = 0.007076820733051469

# Row 7
$ws.Range("I7").Value = 0.003347193013613811
$ws.Range("J7").Value = 0.003347193013613811
$ws.Range("M7").Value = 51.402972
$ws.Range("N7").Value = 154.208916
$ws.Range("O7").Value = 0.2478201393026385
$ws.Range("P7").Value = 0.2478201393026385
$ws.Range("Q7").Value = 18.218721097312
$ws.Range("R7").Value = 163.968489875808
$ws.Range("S7").Value = 0.0008295018389065928
$ws.Range("T7").Value = 0.0008295018389065928

# Row 8
$ws.Range("I8").Value = 0.003347193013613811
$ws.Range("J8").Value = 0.003347193013613811
$ws.Range("O8").Value = 0.2505855219821376
$ws.Range("P8").Value = 0.2505855219821376
$ws.Range("S8").Value = 0.0008387581084913808
$ws.Range("T8").Value = 0.0008387581084913808

# Row 9
$ws.Range("I9").Value = 0.003347193013613811
$ws.Range("J9").Value = 0.003347193013613811
$ws.Range("M9").Value = 59.18509700000001
$ws.Range("N9").Value = 177.555291
$ws.Range("O9").Value = 0.2853387345614992
$ws.Range("P9").Value = 0.2853387345614991
$ws.Range("Q9").Value = 20.97693447297867
$ws.Range("R9").Value = 188.792410256808
$ws.Range("S9").Value = 0.0009550838188376556
$ws.Range("T9").Value = 0.0009550838188376554

# Row 10
$ws.Range("I10").Value = 0.003347193013613811
$ws.Range("J10").Value = 0.003347193013613811
$ws.Range("M10").Value = 8.716382666666666
$ws.Range("N10").Value = 26.149148
$ws.Range("O10").Value = 0.04202276799614693
$ws.Range("P10").Value = 0.04202276799614692
$ws.Range("Q10").Value = 3.089341697624889
$ws.Range("R10").Value = 27.804075278624
$ws.Range("S10").Value = 0.000140658315449417
$ws.Range("T10").Value = 0.000140658315449417

# Row 11
$ws.Range("I11").Value = 0.003347193013613811
$ws.Range("J11").Value = 0.003347193013613811
$ws.Range("M11").Value = 36.13945833333333
$ws.Range("N11").Value = 108.418375
$ws.Range("O11").Value = 0.1742328361575779
$ws.Range("P11").Value = 0.1742328361575779
$ws.Range("Q11").Value = 12.80888412411111
$ws.Range("R11").Value = 115.279957117
$ws.Range("S11").Value = 0.0005831909319287646
$ws.Range("T11").Value = 0.0005831909319287646

# Row 12
$ws.Range("G12").Value = 61.65203333333333
$ws.Range("H12").Value = 184.9561
$ws.Range("I12").Value = 0.582235260574047
$ws.Range("J12").Value = 0.5822352605740471
$ws.Range("M12").Value = 51.402972
$ws.Range("N12").Value = 154.208916
$ws.Range("O12").Value = 0.2478201393026385
$ws.Range("P12").Value = 0.2478201393026385
$ws.Range("Q12").Value = 3169.097743176399
$ws.Range("R12").Value = 28521.8796885876
$ws.Range("S12").Value = 0.1442896233823683
$ws.Range("T12").Value = 0.1442896233823683

# Row 13
$ws.Range("G13").Value = 61.65203333333333
$ws.Range("H13").Value = 184.9561
$ws.Range("I13").Value = 0.582235260574047
$ws.Range("J13").Value = 0.5822352605740471
$ws.Range("O13").Value = 0.2505855219821376
$ws.Range("P13").Value = 0.2505855219821376
$ws.Range("Q13").Value = 3204.4611645403
$ws.Range("R13").Value = 28840.1504808627
$ws.Range("S13").Value = 0.1458997266873534
$ws.Range("T13").Value = 0.1458997266873535

# Row 14
$ws.Range("G14").Value = 61.65203333333333
$ws.Range("H14").Value = 184.9561
$ws.Range("I14").Value = 0.582235260574047
$ws.Range("J14").Value = 0.5822352605740471
$ws.Range("M14").Value = 59.18509700000001
$ws.Range("N14").Value = 177.555291
$ws.Range("O14").Value = 0.2853387345614992
$ws.Range("P14").Value = 0.2853387345614991
$ws.Range("Q14").Value = 3648.881573080567
$ws.Range("R14").Value = 32839.9341577251
$ws.Range("S14").Value = 0.1661342724692833
$ws.Range("T14").Value = 0.1661342724692833

# Row 15
$ws.Range("G15").Value = 61.65203333333333
$ws.Range("H15").Value = 184.9561
$ws.Range("I15").Value = 0.582235260574047
$ws.Range("J15").Value = 0.5822352605740471
$ws.Range("M15").Value = 8.716382666666666
$ws.Range("N15").Value = 26.149148
$ws.Range("O15").Value = 0.04202276799614693
$ws.Range("P15").Value = 0.04202276799614692
$ws.Range("Q15").Value = 537.3827147114222
$ws.Range("R15").Value = 4836.444432402799
$ws.Range("S15").Value = 0.02446713727427933
$ws.Range("T15").Value = 0.02446713727427933

# Row 16
$ws.Range("G16").Value = 61.65203333333333
$ws.Range("H16").Value = 184.9561
$ws.Range("I16").Value = 0.582235260574047
$ws.Range("J16").Value = 0.5822352605740471
$ws.Range("M16").Value = 36.13945833333333
$ws.Range("N16").Value = 108.418375
$ws.Range("O16").Value = 0.1742328361575779
$ws.Range("P16").Value = 0.1742328361575779
$ws.Range("Q16").Value = 2228.071089815277
$ws.Range("R16").Value = 20052.6398083375
$ws.Range("S16").Value = 0.1014445007607626
$ws.Range("T16").Value = 0.1014445007607626

# Row 17
$ws.Range("G17").Value = 0.10468
$ws.Range("H17").Value = 0.31404
$ws.Range("I17").Value = 0.0009885868118471018
$ws.Range("J17").Value = 0.0009885868118471018
$ws.Range("M17").Value = 51.402972
$ws.Range("N17").Value = 154.208916
$ws.Range("O17").Value = 0.2478201393026385
$ws.Range("P17").Value = 0.2478201393026385
$ws.Range("Q17").Value = 5.38086310896
$ws.Range("R17").Value = 48.42776798063999
$ws.Range("S17").Value = 0.0002449917214247
$ws.Range("T17").Value = 0.0002449917214247

# Row 18
$ws.Range("G18").Value = 0.10468
$ws.Range("H18").Value = 0.31404
$ws.Range("I18").Value = 0.0009885868118471018
$ws.Range("J18").Value = 0.0009885868118471018
$ws.Range("O18").Value = 0.2505855219821376
$ws.Range("P18").Value = 0.2505855219821376
$ws.Range("Q18").Value = 5.44090724292
$ws.Range("R18").Value = 48.96816518628
$ws.Range("S18").Value = 0.0002477255422713632
$ws.Range("T18").Value = 0.0002477255422713632

# Row 19
$ws.Range("G19").Value = 0.10468
$ws.Range("H19").Value = 0.31404
$ws.Range("I19").Value = 0.0009885868118471018
$ws.Range("J19").Value = 0.0009885868118471018
$ws.Range("M19").Value = 59.18509700000001
$ws.Range("N19").Value = 177.555291
$ws.Range("O19").Value = 0.2853387345614992
$ws.Range("P19").Value = 0.2853387345614991
$ws.Range("Q19").Value = 6.19549595396
$ws.Range("R19").Value = 55.75946358564
$ws.Range("S19").Value = 0.0002820821098966389
$ws.Range("T19").Value = 0.0002820821098966388

# Row 20
$ws.Range("G20").Value = 0.10468
$ws.Range("H20").Value = 0.31404
$ws.Range("I20").Value = 0.0009885868118471018
$ws.Range("J20").Value = 0.0009885868118471018
$ws.Range("M20").Value = 8.716382666666666
$ws.Range("N20").Value = 26.149148
$ws.Range("O20").Value = 0.04202276799614693
$ws.Range("P20").Value = 0.04202276799614692
$ws.Range("Q20").Value = 0.9124309375466666
$ws.Range("R20").Value = 8.211878437919999
$ws.Range("S20").Value = 0.00004154315423830132
$ws.Range("T20").Value = 0.00004154315423830131

# Row 21
$ws.Range("G21").Value = 0.10468
$ws.Range("H21").Value = 0.31404
$ws.Range("I21").Value = 0.0009885868118471018
$ws.Range("J21").Value = 0.0009885868118471018
$ws.Range("M21").Value = 36.13945833333333
$ws.Range("N21").Value = 108.418375
$ws.Range("O21").Value = 0.1742328361575779
$ws.Range("P21").Value = 0.1742328361575779
$ws.Range("Q21").Value = 3.783078498333333
$ws.Range("R21").Value = 34.047706485
$ws.Range("S21").Value = 0.0001722442840160984
$ws.Range("T21").Value = 0.0001722442840160984

# Row 22
$ws.Range("G22").Value = 39.47650533333334
$ws.Range("H22").Value = 118.429516
$ws.Range("I22").Value = 0.3728119273055513
$ws.Range("J22").Value = 0.3728119273055513
$ws.Range("M22").Value = 51.402972
$ws.Range("N22").Value = 154.208916
$ws.Range("O22").Value = 0.2478201393026385
$ws.Range("P22").Value = 0.2478201393026385
$ws.Range("Q22").Value = 2029.209698307184
$ws.Range("R22").Value = 18262.88728476466
$ws.Range("S22").Value = 0.09239030375854684
$ws.Range("T22").Value = 0.09239030375854684

# Row 23
$ws.Range("G23").Value = 39.47650533333334
$ws.Range("H23").Value = 118.429516
$ws.Range("I23").Value = 0.3728119273055513
$ws.Range("J23").Value = 0.3728119273055513
$ws.Range("O23").Value = 0.2505855219821376
$ws.Range("P23").Value = 0.2505855219821376
$ws.Range("Q23").Value = 2051.853303336869
$ws.Range("R23").Value = 18466.67973003181
$ws.Range("S23").Value = 0.0934212714050283
$ws.Range("T23").Value = 0.0934212714050283

# Row 24
$ws.Range("G24").Value = 39.47650533333334
$ws.Range("H24").Value = 118.429516
$ws.Range("I24").Value = 0.3728119273055513
$ws.Range("J24").Value = 0.3728119273055513
$ws.Range("M24").Value = 59.18509700000001
$ws.Range("N24").Value = 177.555291
$ws.Range("O24").Value = 0.2853387345614992
$ws.Range("P24").Value = 0.2853387345614991
$ws.Range("Q24").Value = 2336.420797374351
$ws.Range("R24").Value = 21027.78717636916
$ws.Range("S24").Value = 0.1063776835667996
$ws.Range("T24").Value = 0.1063776835667996

# Row 25
$ws.Range("G25").Value = 39.47650533333334
$ws.Range("H25").Value = 118.429516
$ws.Range("I25").Value = 0.3728119273055513
$ws.Range("J25").Value = 0.3728119273055513
$ws.Range("M25").Value = 8.716382666666666
$ws.Range("N25").Value = 26.149148
$ws.Range("O25").Value = 0.04202276799614693
$ws.Range("P25").Value = 0.04202276799614692
$ws.Range("Q25").Value = 344.0923268280409
$ws.Range("R25").Value = 3096.830941452368
$ws.Range("S25").Value = 0.01566658912735758
$ws.Range("T25").Value = 0.01566658912735757

# Row 26
$ws.Range("G26").Value = 39.47650533333334
$ws.Range("H26").Value = 118.429516
$ws.Range("I26").Value = 0.3728119273055513
$ws.Range("J26").Value = 0.3728119273055513
$ws.Range("M26").Value = 36.13945833333333
$ws.Range("N26").Value = 108.418375
$ws.Range("O26").Value = 0.1742328361575779
$ws.Range("P26").Value = 0.1742328361575779
$ws.Range("Q26").Value = 1426.659519639611
$ws.Range("R26").Value = 12839.9356767565
$ws.Range("S26").Value = 0.06495607944781898
$ws.Range("T26").Value = 0.06495607944781898
